$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2156.039195262555
$ws.Range("D3").Value = 158.4280988660578

$ws.Range("B4").Value = 903.8118505389275
$ws.Range("D4").Value = 112.6306310535139

$ws.Range("B5").Value = 650.7692219178074
$ws.Range("D5").Value = 0.0001013698630136

$ws.Range("B6").Value = 1409.902009589038
$ws.Range("D6").Value = 41.91803013698635

$ws.Range("B7").Value = 1862.578558904111
$ws.Range("D7").Value = 144.3288136986301

$ws.Range("B8").Value = 2924.263456849315
$ws.Range("D8").Value = 280

$ws.Range("B9").Value = 4096.821336986302
$ws.Range("D9").Value = 280.0000246575343

$ws.Range("F10").Value = 3104696.441178081

$ws.Range("G11").Value = 0.7358194888566576

$ws.Range("F12").Value = 228136.4623671233
$ws.Range("G12").Value = 0.07348108476607029

$ws.Range("G13").Value = 0.1906994263772722
